$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential notice date from 2021-04-30 to 2021-05-03
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-03 for illustrative purposes only and are subject to change."

# Update the numeric figures in columns D and E for rows 2-4
$ws.Range("D2").Value = 0.8443898325850226
$ws.Range("E2").Value = 0.01101011394187679

$ws.Range("D3").Value = 0.1556101674149772
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.9999999999999999
$ws.Range("E4").Value = 0.009296828268123392

# Restore sheet protection (the workbook ships with the sheet protected)
$ws.Protect()
